$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.705.85"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.600.26"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "'0.512"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "'0.0619"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").Value = "'19.55"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "1.825.02"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "1.631.49"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").Value = "'4.05"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "'65.36"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "26.692.67"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "0.0₃0760"
$ws.Range("E18").Value = "  +4.28%  "
$ws.Range("D19").Value = "'7.20"
$ws.Range("E19").Value = "  +3.48%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "'209.40"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").Value = "'142.86"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("D34").Value = "1.293.00"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").Value = "'0.621"
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  +18.50%  "
$ws.Range("D40").Value = "'0.825"
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "'0.784"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "'63.18"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").Value = "1.737.13"
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("D46").Value = "'91.08"
$ws.Range("E46").Value = "  +1.64%  "
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").Value = "'0.101"
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("E51").Value = "  +0.14%  "
